# Lebanon CS template update: Box file paths changed (moved under
# "Box 3EA Team Folder"), a new fs_data_file_path column was added to
# the "path" (formerly "text") sheet, and the active sheet/tab selection
# moved back to that first sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "text" sheet -> rename to "path", add fs_data_file_path column (C),
#    shove model_file_path out to column D, and refresh the Box paths.
# ---------------------------------------------------------------------
$wsPath = $wb.Worksheets.Item("text")
$wsPath.Name = "path"

# Move the existing model_file_path header/value from C to D first.
$wsPath.Range("D1").Value = $wsPath.Range("C1").Value2
$wsPath.Range("C1").Value = "fs_data_file_path"

# Write the new/updated Box paths in the same order the original author
# entered them (header, then the new D2 path, then B2, then C2) so the
# shared-string table comes out in the same order as the source edit.
$wsPath.Range("D2").Value = "/Users/michaelfive/Box/Box 3EA Team Folder/For Zezhen/MR automation/Test Data"
$wsPath.Range("B2").Value = "/Users/michaelfive/Box/Box 3EA Team Folder/3EA Analysis/3EA Lebanon_Analysis/Lebanon_Y1_FA/LBY1_PREIMPUTED_FULL_SPREAD_10-31-2019_mplus.dta"
$wsPath.Range("C2").Value = "/Users/michaelfive/Box/Box 3EA Team Folder/For Zezhen/MR automation/Test Data/CS123_fscores.csv"

# New column C needs the same kind of explicit width the other two data
# columns already carry.
$wsPath.Columns.Item(3).ColumnWidth = 25

# ---------------------------------------------------------------------
# 2. Selection bookkeeping on the other two sheets (content is otherwise
#    unchanged on "subscale" and "model" - only shared-string indices
#    shift, which Excel handles for us automatically).
# ---------------------------------------------------------------------
$wsSubscale = $wb.Worksheets.Item("subscale")
$wsSubscale.Range("D17").Select()

$wsModel = $wb.Worksheets.Item("model")
$wsModel.Range("F2").Select()

# ---------------------------------------------------------------------
# 3. Make "path" the active sheet/tab again, with B2 selected, which
#    also clears tabSelected on "model" and sets the workbook's
#    active tab back to the first sheet.
# ---------------------------------------------------------------------
$wsPath.Range("B2").Select()
$wsPath.Activate()
